# Correção de uma frase: ajusta a redação das mensagens de "abrir" e
# "fechar" catraca na planilha Plan1 (colunas D, linhas 2 e 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("D2").Value = "O cliente deseja abrir a catraca, essa mensagem é concatenada com o id da catraca que quer abrir e é enviada para o servidor principal"
$ws.Range("D3").Value = "O cliente deseja fechar a catraca, essa mensagem é concatenada com o id da catraca que quer fechar e é enviada para o servidor principal"

# A seleção ativa ficou em D4 após a edição.
$ws.Range("D4").Select()
